$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before the existing row 38 (Early Burlat / Segunda, 44519),
# shifting the existing rows 38-50 down to 40-52.
$ws.Rows("38:39").Insert()

# New row 38: Cereza - Lapins - Primera (Region de O'Higgins)
$ws.Cells.Item(38, 1).Value = 4
$ws.Cells.Item(38, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(38, 3).Value = "Los Lagos"
$ws.Cells.Item(38, 4).Value = 44551
$ws.Cells.Item(38, 5).Value = 10
$ws.Cells.Item(38, 6).Value = "Fruta"
$ws.Cells.Item(38, 7).Value = 100103
$ws.Cells.Item(38, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(38, 9).Value = 100103001
$ws.Cells.Item(38, 10).Value = "Cereza"
$ws.Cells.Item(38, 11).Value = "Lapins"
$ws.Cells.Item(38, 12).Value = "Primera"
$ws.Cells.Item(38, 13).Value = 800
$ws.Cells.Item(38, 14).Value = 6500
$ws.Cells.Item(38, 15).Value = 7000
$ws.Cells.Item(38, 16).Value = 6750
$ws.Cells.Item(38, 17).Value = "$/bandeja 10 kilos"
$ws.Cells.Item(38, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(38, 19).Value = 675
$ws.Cells.Item(38, 20).Value = 10

# New row 39: Cereza - Santina - Primera (Region de O'Higgins)
$ws.Cells.Item(39, 1).Value = 4
$ws.Cells.Item(39, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(39, 3).Value = "Los Lagos"
$ws.Cells.Item(39, 4).Value = 44551
$ws.Cells.Item(39, 5).Value = 10
$ws.Cells.Item(39, 6).Value = "Fruta"
$ws.Cells.Item(39, 7).Value = 100103
$ws.Cells.Item(39, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(39, 9).Value = 100103001
$ws.Cells.Item(39, 10).Value = "Cereza"
$ws.Cells.Item(39, 11).Value = "Santina"
$ws.Cells.Item(39, 12).Value = "Primera"
$ws.Cells.Item(39, 13).Value = 800
$ws.Cells.Item(39, 14).Value = 6500
$ws.Cells.Item(39, 15).Value = 7000
$ws.Cells.Item(39, 16).Value = 6750
$ws.Cells.Item(39, 17).Value = "$/bandeja 10 kilos"
$ws.Cells.Item(39, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(39, 19).Value = 675
$ws.Cells.Item(39, 20).Value = 10
